$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("G51").Select()
$r = $win.ScrollRow
Write-Host "ScrollRow now:" $r
